$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.45
$ws.Range("I5").Value = 3.2
$ws.Range("L5").Value = 4
$ws.Range("W5").Value = 6
$ws.Range("AH5").Value = 13
$ws.Range("AJ5").Value = 34
$ws.Range("AY5").Value = 34
$ws.Range("AZ5").Value = 67
$ws.Range("G7").Value = 2.2
$ws.Range("I7").Value = 3.5
$ws.Range("J7").Value = 3.1
$ws.Range("L7").Value = 4.33
$ws.Range("Z7").Value = 21
$ws.Range("AA7").Value = 23
$ws.Range("AG7").Value = 7.5
$ws.Range("AH7").Value = 15
$ws.Range("AI7").Value = 13
$ws.Range("AK7").Value = 34
$ws.Range("AQ7").Value = 51
$ws.Range("AX7").Value = 21
$ws.Range("G10").Value = 3.4
$ws.Range("K10").Value = 1.95
$ws.Range("L10").Value = 3.1
$ws.Range("N10").Value = 8
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 2.63
$ws.Range("Q10").Value = 2.4
$ws.Range("R10").Value = 1.53
$ws.Range("W10").Value = 8
$ws.Range("Y10").Value = 13
$ws.Range("AC10").Value = 7
$ws.Range("AH10").Value = 9.5
$ws.Range("AI10").Value = 10
$ws.Range("AM10").Value = 1250
$ws.Range("AU10").Value = 9
$ws.Range("G11").Value = 2.3
$ws.Range("I11").Value = 3.4
$ws.Range("J11").Value = 3.1
$ws.Range("X11").Value = 10
$ws.Range("AF11").Value = 51
$ws.Range("AH11").Value = 15
$ws.Range("AJ11").Value = 34
$ws.Range("AN11").Value = 4.33
